# Add shading for sinking: introduce a new "SINKING" sheet (a sibling of
# the existing "SLEEP" sheet) with its own start/end timing data, and trim
# the now-unused helper columns (G:I) that were left blank on "SLEEP".

$wb = $excel.ActiveWorkbook

$activity = $wb.Worksheets.Item("ACTIVITY")
$sleep = $wb.Worksheets.Item("SLEEP")

# SLEEP no longer needs the blank G:I helper columns.
$sleep.Range("G2:I9").Clear()

# Create the new SINKING sheet right after SLEEP.
$sinking = $wb.Worksheets.Add([Type]::Missing, $sleep)
$sinking.Name = "SINKING"

# Header row.
$sinking.Range("A1").Value2 = "Sleep"
$sinking.Range("B1").Value2 = "Start"
$sinking.Range("C1").Value2 = "End"
$sinking.Range("D1").Value2 = "Start"
$sinking.Range("E1").Value2 = "End"

# Data rows use the same [m]:ss time format as the SLEEP sheet.
$sinking.Range("A2:E9").NumberFormat = "[m]:ss"

$sinking.Range("A2").Value2 = "Fly 1 "
$sinking.Range("B2").Value2 = 0.1388888888888889
$sinking.Range("C2").Value2 = 0.14583333333333334
$sinking.Range("D2").Value2 = 0.34722222222222227
$sinking.Range("E2").Value2 = 0.38194444444444442

$sinking.Range("A3").Value2 = "Fly 2"
$sinking.Range("B3").Value2 = 0.12313657407407408
$sinking.Range("C3").Value2 = 0.1388888888888889
$sinking.Range("D3").Value2 = 0.38135416666666666
$sinking.Range("E3").Value2 = 0.42372685185185183

$sinking.Range("A4").Value2 = "Fly 3"
$sinking.Range("B4").Value2 = 0.09886574074074074
$sinking.Range("C4").Value2 = 0.1388888888888889
$sinking.Range("D4").Value2 = 0.41666666666666669
$sinking.Range("E4").Value2 = 0.46609953703703705

$sinking.Range("A5").Value2 = "Fly 4"
$sinking.Range("B5").Value2 = 0.09886574074074074
$sinking.Range("C5").Value2 = 0.1388888888888889
$sinking.Range("D5").Value2 = 0.41666666666666669
$sinking.Range("E5").Value2 = 0.46609953703703705

$sinking.Range("A6").Value2 = "Fly 1 "
$sinking.Range("B6").Value2 = 0.1388888888888889
$sinking.Range("C6").Value2 = 0.15606481481481482
$sinking.Range("D6").Value2 = 0.34722222222222227
$sinking.Range("E6").Value2 = 0.35416666666666669

$sinking.Range("A7").Value2 = "Fly 2"
$sinking.Range("B7").Value2 = 0.1388888888888889
$sinking.Range("C7").Value2 = 0.19773148148148148
$sinking.Range("D7").Value2 = 0.41666666666666669
$sinking.Range("E7").Value2 = 0.42372685185185183

$sinking.Range("A8").Value2 = "Fly 3"
$sinking.Range("B8").Value2 = 0.09886574074074074
$sinking.Range("C8").Value2 = 0.1388888888888889
$sinking.Range("D8").Value2 = 0.41666666666666669
$sinking.Range("E8").Value2 = 0.46609953703703705

$sinking.Range("A9").Value2 = "Fly 4"
$sinking.Range("B9").Value2 = 0.09886574074074074
$sinking.Range("C9").Value2 = 0.1388888888888889
$sinking.Range("D9").Value2 = 0.41666666666666669
$sinking.Range("E9").Value2 = 0.46609953703703705

# Restore/update the cursor position (selection) on each sheet.
$null = $activity.Activate()
$null = $activity.Range("N15").Select()

$null = $sleep.Activate()
$null = $sleep.Range("I14").Select()

$null = $sinking.Activate()
$null = $sinking.Range("N20").Select()
